$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = $ws.Range("E2").Text
$ws.Range("F3").Value = $ws.Range("E3").Text
$ws.Range("F4").Value = $ws.Range("E4").Text
